{"js": "// Apply text replacements for the three-digit / one-digit division worksheet.\n// Each pair is a unique, exact-match old->new string (date header + 25 division cells).\nconst replacements = [\n  [\"2026-01-28 Wednesday\", \"2026-01-29 Thursday\"],\n  [\"950\u00f72=475, 0\", \"265\u00f76=44, 1\"],\n  [\"816\u00f75=163, 1\", \"796\u00f73=265, 1\"],\n  [\"145\u00f72=72, 1\", \"660\u00f75=132, 0\"],\n  [\"937\u00f72=468, 1\", \"828\u00f75=165, 3\"],\n  [\"189\u00f72=94, 1\", \"469\u00f73=156, 1\"],\n  [\"243\u00f76=40, 3\", \"465\u00f72=232, 1\"],\n  [\"181\u00f77=25, 6\", \"498\u00f77=71, 1\"],\n  [\"860\u00f78=107, 4\", \"920\u00f75=184, 0\"],\n  [\"202\u00f74=50, 2\", \"705\u00f79=78, 3\"],\n  [\"102\u00f76=17, 0\", \"642\u00f77=91, 5\"],\n  [\"466\u00f74=116, 2\", \"146\u00f72=73, 0\"],\n  [\"637\u00f78=79, 5\", \"602\u00f78=75, 2\"],\n  [\"246\u00f78=30, 6\", \"687\u00f75=137, 2\"],\n  [\"779\u00f78=97, 3\", \"616\u00f73=205, 1\"],\n  [\"910\u00f79=101, 1\", \"731\u00f76=121, 5\"],\n  [\"162\u00f74=40, 2\", \"147\u00f73=49, 0\"],\n  [\"900\u00f76=150, 0\", \"113\u00f78=14, 1\"],\n  [\"742\u00f74=185, 2\", \"948\u00f76=158, 0\"],\n  [\"596\u00f79=66, 2\", \"857\u00f77=122, 3\"],\n  [\"900\u00f77=128, 4\", \"119\u00f75=23, 4\"],\n  [\"132\u00f74=33, 0\", \"532\u00f77=76, 0\"],\n  [\"573\u00f74=143, 1\", \"841\u00f75=168, 1\"],\n  [\"844\u00f74=211, 0\", \"931\u00f76=155, 1\"],\n  [\"668\u00f76=111, 2\", \"140\u00f76=23, 2\"],\n  [\"664\u00f75=132, 4\", \"153\u00f73=51, 0\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const result of results.items) {\n    result.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update master worksheet: header date + all 25 three-digit / one-digit division problems.\n# Each row is an exact, unique old->new text pair pulled from the target OOXML diff.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"2026-01-28 Wednesday\", \"2026-01-29 Thursday\"),\n  @(\"950\u00f72=475, 0\", \"265\u00f76=44, 1\"),\n  @(\"816\u00f75=163, 1\", \"796\u00f73=265, 1\"),\n  @(\"145\u00f72=72, 1\", \"660\u00f75=132, 0\"),\n  @(\"937\u00f72=468, 1\", \"828\u00f75=165, 3\"),\n  @(\"189\u00f72=94, 1\", \"469\u00f73=156, 1\"),\n  @(\"243\u00f76=40, 3\", \"465\u00f72=232, 1\"),\n  @(\"181\u00f77=25, 6\", \"498\u00f77=71, 1\"),\n  @(\"860\u00f78=107, 4\", \"920\u00f75=184, 0\"),\n  @(\"202\u00f74=50, 2\", \"705\u00f79=78, 3\"),\n  @(\"102\u00f76=17, 0\", \"642\u00f77=91, 5\"),\n  @(\"466\u00f74=116, 2\", \"146\u00f72=73, 0\"),\n  @(\"637\u00f78=79, 5\", \"602\u00f78=75, 2\"),\n  @(\"246\u00f78=30, 6\", \"687\u00f75=137, 2\"),\n  @(\"779\u00f78=97, 3\", \"616\u00f73=205, 1\"),\n  @(\"910\u00f79=101, 1\", \"731\u00f76=121, 5\"),\n  @(\"162\u00f74=40, 2\", \"147\u00f73=49, 0\"),\n  @(\"900\u00f76=150, 0\", \"113\u00f78=14, 1\"),\n  @(\"742\u00f74=185, 2\", \"948\u00f76=158, 0\"),\n  @(\"596\u00f79=66, 2\", \"857\u00f77=122, 3\"),\n  @(\"900\u00f77=128, 4\", \"119\u00f75=23, 4\"),\n  @(\"132\u00f74=33, 0\", \"532\u00f77=76, 0\"),\n  @(\"573\u00f74=143, 1\", \"841\u00f75=168, 1\"),\n  @(\"844\u00f74=211, 0\", \"931\u00f76=155, 1\"),\n  @(\"668\u00f76=111, 2\", \"140\u00f76=23, 2\"),\n  @(\"664\u00f75=132, 4\", \"153\u00f73=51, 0\"),\n)\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
